# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column G holds the "K" values (header at G1 == "K").
# Update the recalculated K values for each data row (row 7 already matches and is left untouched).
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 0
